$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old "Zero" and "ZeroAcc" rows (rows 4 and 5) entirely.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Re-add the example signal rows starting at row 5: Zero, Triangle, Meander,
# Spike.
$data = @(
    @("Zero",     0, 2,    3, 2,    0, 40, 0, 1, 3),
    @("Triangle", 0, 3,    0, 2,    0, 10, 0, 1, 0),
    @("Meander",  1, 0.02, 3, 0.02, 0, 10, 0, 1, 1),
    @("Spike",    2, 0.02, 0, 0.02, 0, 10, 0, 1, 2)
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
}

$ws.Range("F15").Select()
